# Apply the commit's data updates to the "展览" (Exhibitions) and
# "全部类型" (All Types) sheets: refresh "want-to-go" counts for several
# existing events, and insert a new 吉安·COMIC LIFE event row ahead of the
# existing 南昌·萌卡动漫展 / 江西·JMG rows (which shift down by one row and
# get their own counts bumped too).

$wb = $excel.ActiveWorkbook

# ============ Sheet "展览" ============
$ws1 = $wb.Worksheets.Item("展览")

# Refresh "want to go" counts (column F) for existing rows
$ws1.Cells.Item(3,6).Value = 1010
$ws1.Cells.Item(4,6).Value = 163
$ws1.Cells.Item(5,6).Value = 2720
$ws1.Cells.Item(6,6).Value = 90
$ws1.Cells.Item(7,6).Value = 209
$ws1.Cells.Item(9,6).Value = 111

# Insert a new row before row 11 for the new 吉安 event; this pushes the
# old row 11 (南昌·萌卡动漫展) down to row 12 and the old row 12
# (江西·JMG) down to row 13, preserving their data.
$ws1.Rows.Item(11).Insert()

# Give the new row 11 the same index-column formatting as the other rows
$ws1.Cells.Item(11,1).Value = 10
$ws1.Cells.Item(11,1).Font.Bold = $true
$ws1.Cells.Item(11,1).HorizontalAlignment = -4108
$ws1.Cells.Item(11,1).VerticalAlignment = -4160
$ws1.Cells.Item(11,1).Borders.LineStyle = 1

# Fill the newly inserted row 11 with the new event data
$ws1.Range("B11").NumberFormat = "@"
$ws1.Range("B11").Value = "2024-10-01"
$ws1.Range("B11").Style = "Normal"
$ws1.Range("C11").Value = "吉安·COMIC LIFE次元假日06"
$ws1.Range("D11").Value = "东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心"
$ws1.Range("E11").Value = "2024.10.01 09:00-10.01 18:00"
$ws1.Range("F11").Value = 0
$ws1.Range("G11").Value = 9.9
$ws1.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=90901"
$ws1.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202408/J892vhG21724034395965.jpeg"

# Update the "want to go" counts for the two rows that shifted down
$ws1.Cells.Item(12,6).Value = 2538   # 南昌·萌卡动漫展 (was row 11)
$ws1.Cells.Item(13,6).Value = 656    # 江西·JMG (was row 12)

# ============ Sheet "全部类型" ============
$ws4 = $wb.Worksheets.Item("全部类型")

# Refresh "want to go" counts (column F) for existing rows
$ws4.Cells.Item(4,6).Value = 1010
$ws4.Cells.Item(5,6).Value = 163
$ws4.Cells.Item(6,6).Value = 2720
$ws4.Cells.Item(7,6).Value = 90
$ws4.Cells.Item(8,6).Value = 209
$ws4.Cells.Item(11,6).Value = 111

# Insert a new row before row 13 for the new 吉安 event; this pushes the
# old row 13 (南昌·萌卡动漫展) down to row 14 and the old row 14
# (江西·JMG) down to row 15, preserving their data.
$ws4.Rows.Item(13).Insert()

# Give the new row 13 the same index-column formatting as the other rows
$ws4.Cells.Item(13,1).Value = 12
$ws4.Cells.Item(13,1).Font.Bold = $true
$ws4.Cells.Item(13,1).HorizontalAlignment = -4108
$ws4.Cells.Item(13,1).VerticalAlignment = -4160
$ws4.Cells.Item(13,1).Borders.LineStyle = 1

# Fill the newly inserted row 13 with the new event data
$ws4.Range("B13").NumberFormat = "@"
$ws4.Range("B13").Value = "2024-10-01"
$ws4.Range("B13").Style = "Normal"
$ws4.Range("C13").Value = "吉安·COMIC LIFE次元假日06"
$ws4.Range("D13").Value = "东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心"
$ws4.Range("E13").Value = "2024.10.01 09:00-10.01 18:00"
$ws4.Range("F13").Value = 0
$ws4.Range("G13").Value = 9.9
$ws4.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=90901"
$ws4.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202408/J892vhG21724034395965.jpeg"

# Update the "want to go" counts for the two rows that shifted down
$ws4.Cells.Item(14,6).Value = 2538   # 南昌·萌卡动漫展 (was row 13)
$ws4.Cells.Item(15,6).Value = 656    # 江西·JMG (was row 14)

